$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric metric values for random_forest, lsboost, neural_network rows.
# (Row 5 "old_model" values stay the same - only its label cell's shared string id
#  changes in the underlying XML, which has no visible effect.)

# Row 2: random_forest
$ws.Range("B2").Value = 4.7087571496212126
$ws.Range("C2").Value = 0.49917528757569113
$ws.Range("D2").Value = 3.7736643984630578
$ws.Range("E2").Value = 0.45409481108355809
$ws.Range("F2").Value = 0.67386557345182585
$ws.Range("G2").Value = 0.60702192714933068
$ws.Range("H2").Value = 0.54590518891644191
$ws.Range("I2").Value = 0.7928388375009936

# Row 3: lsboost
$ws.Range("B3").Value = 4.9708505991732981
$ws.Range("C3").Value = 0.52695981094241251
$ws.Range("D3").Value = 3.8836863321307282
$ws.Range("E3").Value = 0.50605226719204566
$ws.Range("F3").Value = 0.71137350751349016
$ws.Range("G3").Value = 0.62471977177770988
$ws.Range("H3").Value = 0.49394773280795434
$ws.Range("I3").Value = 0.74325595572698289

# Row 4: neural_network
$ws.Range("B4").Value = 4.2753486711292013
$ws.Range("C4").Value = 0.45322966009596516
$ws.Range("D4").Value = 3.2670238705032988
$ws.Range("E4").Value = 0.37434930553829177
$ws.Range("F4").Value = 0.6118409152208536
$ws.Range("G4").Value = 0.52552503787127436
$ws.Range("H4").Value = 0.62565069446170818
$ws.Range("I4").Value = 0.80887238697455244

$wb.Save()
